# Applies the "Add files via upload" change to the Saldo export sheet:
#   1. Insert a new account row (004241147 / ANTONIO / 100000) right
#      before the existing "008054713 / MODULAR" row (Excel row 4).
#   2. Insert a new account row (001368670 / THIAGO / 97.2) right before
#      the existing "004239387 / LUIZ" row (originally Excel row 27, now
#      row 28 after step 1 shifted everything down by one).
#   3. Remove the old "001368670 / THIAGO / -2.8" row that used to sit
#      near the bottom of the sheet, just above the trailing blank row
#      and the "Filtros aplicados" footer row (originally Excel row 234,
#      now row 236 after the two inserts above).
#
# Net effect: row count goes from 236 -> 237 (+1 +1 -1), and the THIAGO
# account ends up with a single, relocated row (97.2) instead of its
# old entry (-2.8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert ANTONIO row before "008054713" (Excel row 4) ----------------
$ws.Rows.Item(4).EntireRow.Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"   # keep the leading zeros as text
$ws.Cells.Item(4, 1).Value = "004241147"
$ws.Cells.Item(4, 2).Value = "ANTONIO"
$ws.Cells.Item(4, 3).Value = 100000

# --- 2) Insert THIAGO row before "004239387" (now Excel row 28) ------------
$ws.Rows.Item(28).EntireRow.Insert()
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "001368670"
$ws.Cells.Item(28, 2).Value = "THIAGO"
$ws.Cells.Item(28, 3).Value = 97.2

# --- 3) Delete the old THIAGO (-2.8) row (now Excel row 236) ---------------
$ws.Rows.Item(236).EntireRow.Delete()
